$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF). Copy the formatting from the
# existing header cell (H1) so the new headers share the bold/centered/
# bordered header style, then set their text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new data columns I (I0) and J (IF) for rows 2-20.
$data = @(
    @(2, 8, 8),
    @(3, 8, 9),
    @(4, 7, 8),
    @(5, 7, 7),
    @(6, 7, 7),
    @(7, 6, 7),
    @(8, 5, 5),
    @(9, 6, 6),
    @(10, 9, 9),
    @(11, 8, 8),
    @(12, 8, 8),
    @(13, 5, 5),
    @(14, 7, 8),
    @(15, 8, 8),
    @(16, 4, 5),
    @(17, 6, 6),
    @(18, 9, 9),
    @(19, 8, 8),
    @(20, 7, 7)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
